# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'66.915.52"
$ws.Cells.Item(2, 5).Value = "  +3.14%  "
$ws.Cells.Item(3, 4).Value = "'3.439.19"
$ws.Cells.Item(3, 5).Value = "  +1.65%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'570.45"
$ws.Cells.Item(5, 5).Value = "  +2.07%  "
$ws.Cells.Item(6, 4).Value = "'184.82"
$ws.Cells.Item(6, 5).Value = "  +5.44%  "
$ws.Cells.Item(7, 5).Value = "  +1.71%  "
$ws.Cells.Item(8, 4).Value = "'3.434.00"
$ws.Cells.Item(8, 5).Value = "  +1.71%  "
$ws.Cells.Item(9, 5).Value = "  +0.03%  "
$ws.Cells.Item(10, 4).Value = "'0.177"
$ws.Cells.Item(10, 5).Value = "  +7.09%  "
$ws.Cells.Item(11, 5).Value = "  +1.79%  "
$ws.Cells.Item(12, 4).Value = "'55.36"
$ws.Cells.Item(12, 5).Value = "  +2.27%  "
$ws.Cells.Item(13, 5).Value = "  +1.78%  "
$ws.Cells.Item(14, 5).Value = "  +3.03%  "
$ws.Cells.Item(15, 4).Value = "'3.988.27"
$ws.Cells.Item(15, 5).Value = "  +1.68%  "
$ws.Cells.Item(16, 4).Value = "'18.54"
$ws.Cells.Item(16, 5).Value = "  +1.46%  "
$ws.Cells.Item(17, 4).Value = "'3.440.71"
$ws.Cells.Item(17, 5).Value = "  +2.13%  "
$ws.Cells.Item(18, 5).Value = "  +0.46%  "
$ws.Cells.Item(19, 4).Value = "'66.756.96"
$ws.Cells.Item(19, 5).Value = "  +3.01%  "
$ws.Cells.Item(20, 4).Value = "'12.00"
$ws.Cells.Item(20, 5).Value = "  +1.78%  "
$ws.Cells.Item(21, 5).Value = "  +1.86%  "
$ws.Cells.Item(22, 4).Value = "'476.54"
$ws.Cells.Item(22, 5).Value = "  +3.28%  "
$ws.Cells.Item(23, 4).Value = "'4.97"
$ws.Cells.Item(23, 5).Value = "  +2.14%  "
$ws.Cells.Item(24, 4).Value = "'14.94"
$ws.Cells.Item(24, 5).Value = "  +10.69%  "
$ws.Cells.Item(25, 4).Value = "'4.20"
$ws.Cells.Item(25, 5).Value = "  +1.74%  "
$ws.Cells.Item(26, 4).Value = "'89.56"
$ws.Cells.Item(26, 5).Value = "  +3.55%  "
$ws.Cells.Item(27, 4).Value = "'2.95"
$ws.Cells.Item(27, 5).Value = "  +0.11%  "
$ws.Cells.Item(28, 4).Value = "'10.99"
$ws.Cells.Item(28, 5).Value = "  +1.42%  "
$ws.Cells.Item(29, 5).Value = "  +2.19%  "
$ws.Cells.Item(30, 4).Value = "'31.52"
$ws.Cells.Item(30, 5).Value = "  +2.37%  "
$ws.Cells.Item(31, 4).Value = "'6.95"
$ws.Cells.Item(31, 5).Value = "  +3.17%  "
$ws.Cells.Item(32, 5).Value = "  +1.29%  "
$ws.Cells.Item(33, 4).Value = "'589.79"
$ws.Cells.Item(33, 5).Value = "  +3.45%  "
$ws.Cells.Item(34, 4).Value = "'63.13"
$ws.Cells.Item(34, 5).Value = "  +3.07%  "
$ws.Cells.Item(35, 5).Value = "  +1.52%  "
$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).Value = "'0.148"
$ws.Cells.Item(36, 5).Value = "  +5.89%  "
$ws.Cells.Item(37, 2).Value = "Dai"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(37, 4).Value = "'0.999"
$ws.Cells.Item(37, 5).Value = "  -0.08%  "
$ws.Cells.Item(38, 4).Value = "'3.64"
$ws.Cells.Item(38, 5).Value = "  -0.02%  "
$ws.Cells.Item(39, 4).Value = "'0.391"
$ws.Cells.Item(39, 5).Value = "  +5.75%  "
$ws.Cells.Item(40, 5).Value = "  +3.06%  "
$ws.Cells.Item(41, 4).Value = "'0.0₃0773"
$ws.Cells.Item(41, 5).Value = "  +3.88%  "
$ws.Cells.Item(42, 4).Value = "'3.128.22"
$ws.Cells.Item(42, 5).Value = "  +1.59%  "
$ws.Cells.Item(43, 4).Value = "'2.92"
$ws.Cells.Item(43, 5).Value = "  +2.85%  "
$ws.Cells.Item(44, 5).Value = "  +7.50%  "
$ws.Cells.Item(45, 4).Value = "'0.0424"
$ws.Cells.Item(45, 5).Value = "  +2.21%  "
$ws.Cells.Item(46, 4).Value = "'2.79"
$ws.Cells.Item(46, 5).Value = "  +19.73%  "
$ws.Cells.Item(47, 4).Value = "'3.26"
$ws.Cells.Item(47, 5).Value = "  +4.56%  "
$ws.Cells.Item(48, 5).Value = "  +0.38%  "
$ws.Cells.Item(49, 5).Value = "  -0.03%  "
$ws.Cells.Item(50, 4).Value = "'141.94"
$ws.Cells.Item(50, 5).Value = "  +2.87%  "
$ws.Cells.Item(51, 5).Value = "  +5.52%  "

Write-Host "Update complete"
